# Adds the 2024-06-09 incremental-extract data for the violent-crime-full-year
# workbook: the 2024 year-to-date column (K, and occasionally J for the 2023
# column where late-reported incidents shifted a prior year's tally) on the
# 'Citywide Totals' / 'By Neighborhood' summary sheets and every affected
# per-neighborhood sheet. Category rows 2-6 hold the per-crime-type counts and
# row 7 (row 101 on 'By Neighborhood') holds the Total.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 11).Value = 3333  # K2: 3304 -> 3333
$ws.Cells.Item(3, 10).Value = 8080  # J3: 8079 -> 8080
$ws.Cells.Item(3, 11).Value = 3296  # K3: 3277 -> 3296
$ws.Cells.Item(4, 11).Value = 684  # K4: 679 -> 684
$ws.Cells.Item(5, 11).Value = 214  # K5: 213 -> 214
$ws.Cells.Item(6, 11).Value = 3890  # K6: 3864 -> 3890
$ws.Cells.Item(7, 10).Value = 29288  # J7: 29287 -> 29288
$ws.Cells.Item(7, 11).Value = 11417  # K7: 11337 -> 11417

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(2, 11).Value = 92  # K2: 91 -> 92
$ws.Cells.Item(4, 11).Value = 39  # K4: 38 -> 39
$ws.Cells.Item(7, 11).Value = 324  # K7: 319 -> 324
$ws.Cells.Item(8, 11).Value = 757  # K8: 751 -> 757
$ws.Cells.Item(10, 11).Value = 64  # K10: 63 -> 64
$ws.Cells.Item(11, 11).Value = 242  # K11: 240 -> 242
$ws.Cells.Item(14, 11).Value = 58  # K14: 57 -> 58
$ws.Cells.Item(15, 11).Value = 117  # K15: 116 -> 117
$ws.Cells.Item(16, 11).Value = 39  # K16: 38 -> 39
$ws.Cells.Item(18, 11).Value = 78  # K18: 77 -> 78
$ws.Cells.Item(20, 11).Value = 261  # K20: 260 -> 261
$ws.Cells.Item(23, 11).Value = 110  # K23: 109 -> 110
$ws.Cells.Item(24, 11).Value = 37  # K24: 36 -> 37
$ws.Cells.Item(25, 11).Value = 47  # K25: 46 -> 47
$ws.Cells.Item(29, 11).Value = 597  # K29: 591 -> 597
$ws.Cells.Item(31, 11).Value = 124  # K31: 123 -> 124
$ws.Cells.Item(33, 11).Value = 452  # K33: 449 -> 452
$ws.Cells.Item(37, 11).Value = 398  # K37: 392 -> 398
$ws.Cells.Item(40, 11).Value = 29  # K40: 28 -> 29
$ws.Cells.Item(42, 11).Value = 408  # K42: 404 -> 408
$ws.Cells.Item(43, 11).Value = 102  # K43: 100 -> 102
$ws.Cells.Item(47, 11).Value = 62  # K47: 60 -> 62
$ws.Cells.Item(51, 11).Value = 135  # K51: 132 -> 135
$ws.Cells.Item(52, 11).Value = 310  # K52: 306 -> 310
$ws.Cells.Item(54, 11).Value = 222  # K54: 221 -> 222
$ws.Cells.Item(56, 11).Value = 13  # K56: 12 -> 13
$ws.Cells.Item(63, 11).Value = 35  # K63: 38 -> 35
$ws.Cells.Item(64, 11).Value = 72  # K64: 69 -> 72
$ws.Cells.Item(65, 11).Value = 272  # K65: 270 -> 272
$ws.Cells.Item(70, 11).Value = 21  # K70: 19 -> 21
$ws.Cells.Item(76, 11).Value = 173  # K76: 172 -> 173
$ws.Cells.Item(79, 11).Value = 296  # K79: 295 -> 296
$ws.Cells.Item(83, 11).Value = 247  # K83: 246 -> 247
$ws.Cells.Item(84, 11).Value = 80  # K84: 77 -> 80
$ws.Cells.Item(85, 11).Value = 541  # K85: 539 -> 541
$ws.Cells.Item(86, 11).Value = 76  # K86: 75 -> 76
$ws.Cells.Item(88, 11).Value = 136  # K88: 135 -> 136
$ws.Cells.Item(89, 11).Value = 151  # K89: 150 -> 151
$ws.Cells.Item(90, 11).Value = 101  # K90: 100 -> 101
$ws.Cells.Item(94, 11).Value = 141  # K94: 140 -> 141
$ws.Cells.Item(96, 11).Value = 140  # K96: 139 -> 140
$ws.Cells.Item(97, 10).Value = 261  # J97: 260 -> 261
$ws.Cells.Item(97, 11).Value = 100  # K97: 99 -> 100
$ws.Cells.Item(98, 11).Value = 62  # K98: 61 -> 62
$ws.Cells.Item(99, 11).Value = 199  # K99: 197 -> 199
$ws.Cells.Item(101, 10).Value = 29288  # J101: 29287 -> 29288
$ws.Cells.Item(101, 11).Value = 11417  # K101: 11337 -> 11417

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Cells.Item(6, 11).Value = 19  # K6: 18 -> 19
$ws.Cells.Item(7, 11).Value = 58  # K7: 57 -> 58

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Cells.Item(6, 11).Value = 68  # K6: 67 -> 68
$ws.Cells.Item(7, 11).Value = 140  # K7: 139 -> 140

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(2, 11).Value = 116  # K2: 112 -> 116
$ws.Cells.Item(6, 11).Value = 84  # K6: 83 -> 84
$ws.Cells.Item(7, 11).Value = 324  # K7: 319 -> 324

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Cells.Item(2, 11).Value = 72  # K2: 71 -> 72
$ws.Cells.Item(4, 11).Value = 12  # K4: 11 -> 12
$ws.Cells.Item(7, 11).Value = 242  # K7: 240 -> 242

$ws = $wb.Worksheets.Item('Uptown')
$ws.Cells.Item(3, 11).Value = 49  # K3: 48 -> 49
$ws.Cells.Item(7, 11).Value = 151  # K7: 150 -> 151

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(2, 11).Value = 192  # K2: 191 -> 192
$ws.Cells.Item(6, 11).Value = 125  # K6: 124 -> 125
$ws.Cells.Item(7, 11).Value = 541  # K7: 539 -> 541

$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(2, 11).Value = 81  # K2: 80 -> 81
$ws.Cells.Item(4, 11).Value = 15  # K4: 14 -> 15
$ws.Cells.Item(6, 11).Value = 124  # K6: 122 -> 124
$ws.Cells.Item(7, 11).Value = 310  # K7: 306 -> 310

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(2, 11).Value = 220  # K2: 218 -> 220
$ws.Cells.Item(3, 11).Value = 230  # K3: 228 -> 230
$ws.Cells.Item(6, 11).Value = 245  # K6: 243 -> 245
$ws.Cells.Item(7, 11).Value = 757  # K7: 751 -> 757

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(3, 11).Value = 80  # K3: 79 -> 80
$ws.Cells.Item(7, 11).Value = 247  # K7: 246 -> 247

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(3, 11).Value = 165  # K3: 163 -> 165
$ws.Cells.Item(6, 11).Value = 130  # K6: 129 -> 130
$ws.Cells.Item(7, 11).Value = 452  # K7: 449 -> 452

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(2, 11).Value = 108  # K2: 107 -> 108
$ws.Cells.Item(3, 11).Value = 137  # K3: 135 -> 137
$ws.Cells.Item(6, 11).Value = 122  # K6: 119 -> 122
$ws.Cells.Item(7, 11).Value = 398  # K7: 392 -> 398

$ws = $wb.Worksheets.Item('New City')
$ws.Cells.Item(6, 11).Value = 106  # K6: 104 -> 106
$ws.Cells.Item(7, 11).Value = 272  # K7: 270 -> 272

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(3, 11).Value = 74  # K3: 73 -> 74
$ws.Cells.Item(5, 11).Value = 6  # K5: 5 -> 6
$ws.Cells.Item(7, 11).Value = 199  # K7: 197 -> 199

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Cells.Item(2, 11).Value = 42  # K2: 41 -> 42
$ws.Cells.Item(7, 11).Value = 124  # K7: 123 -> 124

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(4, 11).Value = 24  # K4: 25 -> 24
$ws.Cells.Item(6, 11).Value = 129  # K6: 128 -> 129

$ws = $wb.Worksheets.Item('South Deering')
$ws.Cells.Item(2, 11).Value = 25  # K2: 24 -> 25
$ws.Cells.Item(6, 11).Value = 20  # K6: 18 -> 20
$ws.Cells.Item(7, 11).Value = 80  # K7: 77 -> 80

$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(3, 11).Value = 66  # K3: 65 -> 66
$ws.Cells.Item(7, 11).Value = 222  # K7: 221 -> 222

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(2, 11).Value = 167  # K2: 163 -> 167
$ws.Cells.Item(6, 11).Value = 181  # K6: 179 -> 181
$ws.Cells.Item(7, 11).Value = 597  # K7: 591 -> 597

$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(6, 11).Value = 102  # K6: 101 -> 102
$ws.Cells.Item(7, 11).Value = 173  # K7: 172 -> 173

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(2, 11).Value = 107  # K2: 104 -> 107
$ws.Cells.Item(3, 11).Value = 130  # K3: 129 -> 130
$ws.Cells.Item(7, 11).Value = 408  # K7: 404 -> 408

$ws = $wb.Worksheets.Item('Avondale')
$ws.Cells.Item(2, 11).Value = 21  # K2: 20 -> 21
$ws.Cells.Item(7, 11).Value = 64  # K7: 63 -> 64

$ws = $wb.Worksheets.Item('Dunning')
$ws.Cells.Item(2, 11).Value = 13  # K2: 12 -> 13
$ws.Cells.Item(7, 11).Value = 37  # K7: 36 -> 37

$ws = $wb.Worksheets.Item('Douglas')
$ws.Cells.Item(3, 11).Value = 37  # K3: 36 -> 37
$ws.Cells.Item(7, 11).Value = 110  # K7: 109 -> 110

$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(3, 11).Value = 102  # K3: 101 -> 102
$ws.Cells.Item(7, 11).Value = 296  # K7: 295 -> 296

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Cells.Item(3, 11).Value = 23  # K3: 22 -> 23
$ws.Cells.Item(6, 11).Value = 28  # K6: 26 -> 28
$ws.Cells.Item(7, 11).Value = 72  # K7: 69 -> 72

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(2, 11).Value = 91  # K2: 90 -> 91
$ws.Cells.Item(7, 11).Value = 261  # K7: 260 -> 261

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Cells.Item(4, 11).Value = 10  # K4: 9 -> 10
$ws.Cells.Item(7, 11).Value = 78  # K7: 77 -> 78

$ws = $wb.Worksheets.Item('West Loop')
$ws.Cells.Item(4, 11).Value = 14  # K4: 13 -> 14
$ws.Cells.Item(7, 11).Value = 141  # K7: 140 -> 141

$ws = $wb.Worksheets.Item('East Side')
$ws.Cells.Item(3, 11).Value = 19  # K3: 18 -> 19
$ws.Cells.Item(7, 11).Value = 47  # K7: 46 -> 47

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Cells.Item(2, 11).Value = 18  # K2: 16 -> 18
$ws.Cells.Item(7, 11).Value = 62  # K7: 60 -> 62

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Cells.Item(3, 11).Value = 30  # K3: 29 -> 30
$ws.Cells.Item(7, 11).Value = 117  # K7: 116 -> 117

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Cells.Item(3, 11).Value = 9  # K3: 8 -> 9
$ws.Cells.Item(7, 11).Value = 62  # K7: 61 -> 62

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Cells.Item(3, 11).Value = 25  # K3: 24 -> 25
$ws.Cells.Item(7, 11).Value = 92  # K7: 91 -> 92

$ws = $wb.Worksheets.Item('West Town')
$ws.Cells.Item(3, 10).Value = 26  # J3: 25 -> 26
$ws.Cells.Item(3, 11).Value = 15  # K3: 14 -> 15
$ws.Cells.Item(7, 10).Value = 261  # J7: 260 -> 261
$ws.Cells.Item(7, 11).Value = 100  # K7: 99 -> 100

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Cells.Item(2, 11).Value = 8  # K2: 7 -> 8
$ws.Cells.Item(4, 11).Value = 6  # K4: 5 -> 6
$ws.Cells.Item(7, 11).Value = 21  # K7: 19 -> 21

$ws = $wb.Worksheets.Item('United Center')
$ws.Cells.Item(2, 11).Value = 32  # K2: 30 -> 32
$ws.Cells.Item(3, 11).Value = 35  # K3: 36 -> 35
$ws.Cells.Item(7, 11).Value = 136  # K7: 135 -> 136

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Cells.Item(6, 11).Value = 19  # K6: 18 -> 19
$ws.Cells.Item(7, 11).Value = 76  # K7: 75 -> 76

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Cells.Item(4, 11).Value = 9  # K4: 8 -> 9
$ws.Cells.Item(7, 11).Value = 101  # K7: 100 -> 101

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Cells.Item(3, 11).Value = 34  # K3: 32 -> 34
$ws.Cells.Item(6, 11).Value = 49  # K6: 48 -> 49
$ws.Cells.Item(7, 11).Value = 135  # K7: 132 -> 135

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Cells.Item(2, 11).Value = 8  # K2: 7 -> 8
$ws.Cells.Item(6, 11).Value = 20  # K6: 21 -> 20

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Cells.Item(4, 11).Value = 11  # K4: 10 -> 11
$ws.Cells.Item(6, 11).Value = 43  # K6: 42 -> 43
$ws.Cells.Item(7, 11).Value = 102  # K7: 100 -> 102

$ws = $wb.Worksheets.Item('Magnificent Mile')
$ws.Cells.Item(2, 11).Value = 1  # K2: None -> 1
$ws.Cells.Item(7, 11).Value = 13  # K7: 12 -> 13

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Cells.Item(6, 11).Value = 6  # K6: 5 -> 6
$ws.Cells.Item(7, 11).Value = 29  # K7: 28 -> 29

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Cells.Item(6, 11).Value = 15  # K6: 14 -> 15
$ws.Cells.Item(7, 11).Value = 39  # K7: 38 -> 39

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Cells.Item(3, 11).Value = 2  # K3: 1 -> 2
$ws.Cells.Item(7, 11).Value = 39  # K7: 38 -> 39

